# =====================================================================================
# This script applies a scraper re-run update to the Serie A 2023-2024 odds workbook:
#  1) Six pairs/triples of rows had their fixture data (columns F:V) rotated because the
#     earlier scrape had them assigned to the wrong row index (date/index columns A:E
#     were already correct and are left untouched).
#  2) Ten new matches (rows 92-101, Indice 91-100) scraped on 31-10-2023 are appended.
# =====================================================================================

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---- Part 1: fix rows whose F:V data block belongs to a different row in the group ----

# Rows 8, 9
$ws.Cells.Item(8,6).Value = "Udinese"
$ws.Cells.Item(8,7).Value = 0
$ws.Cells.Item(8,8).Value = "Juventus"
$ws.Cells.Item(8,9).Value = 3
$ws.Cells.Item(8,10).Value = 4.33
$ws.Cells.Item(8,11).Value = "05/07/2023 21:25"
$ws.Cells.Item(8,12).Value = 4.05
$ws.Cells.Item(8,13).Value = "20/08/2023 20:40"
$ws.Cells.Item(8,14).Value = 3.46
$ws.Cells.Item(8,15).Value = "05/07/2023 21:25"
$ws.Cells.Item(8,16).Value = 3.41
$ws.Cells.Item(8,17).Value = "20/08/2023 20:31"
$ws.Cells.Item(8,18).Value = 1.85
$ws.Cells.Item(8,19).Value = "05/07/2023 21:25"
$ws.Cells.Item(8,20).Value = 2.07
$ws.Cells.Item(8,21).Value = "20/08/2023 20:31"
$ws.Cells.Item(8,22).Value = "https://www.betexplorer.com/football/italy/serie-a/udinese-juventus/zFCpUwv1/"
$ws.Cells.Item(9,6).Value = "Lecce"
$ws.Cells.Item(9,7).Value = 2
$ws.Cells.Item(9,8).Value = "Lazio"
$ws.Cells.Item(9,9).Value = 1
$ws.Cells.Item(9,10).Value = 3.66
$ws.Cells.Item(9,11).Value = "05/07/2023 21:25"
$ws.Cells.Item(9,12).Value = 3.76
$ws.Cells.Item(9,13).Value = "20/08/2023 20:44"
$ws.Cells.Item(9,14).Value = 3.11
$ws.Cells.Item(9,15).Value = "05/07/2023 21:25"
$ws.Cells.Item(9,16).Value = 3.34
$ws.Cells.Item(9,17).Value = "20/08/2023 20:43"
$ws.Cells.Item(9,18).Value = 2.28
$ws.Cells.Item(9,19).Value = "05/07/2023 21:25"
$ws.Cells.Item(9,20).Value = 2.19
$ws.Cells.Item(9,21).Value = "20/08/2023 20:42"
$ws.Cells.Item(9,22).Value = "https://www.betexplorer.com/football/italy/serie-a/lecce-lazio/80d4l8PF/"

# Rows 24, 25
$ws.Cells.Item(24,6).Value = "Bologna"
$ws.Cells.Item(24,7).Value = 2
$ws.Cells.Item(24,8).Value = "Cagliari"
$ws.Cells.Item(24,9).Value = 1
$ws.Cells.Item(24,10).Value = 1.85
$ws.Cells.Item(24,11).Value = "22/08/2023 14:46"
$ws.Cells.Item(24,12).Value = 1.81
$ws.Cells.Item(24,13).Value = "02/09/2023 18:28"
$ws.Cells.Item(24,14).Value = 3.33
$ws.Cells.Item(24,15).Value = "22/08/2023 14:46"
$ws.Cells.Item(24,16).Value = 3.88
$ws.Cells.Item(24,17).Value = "02/09/2023 18:29"
$ws.Cells.Item(24,18).Value = 4.56
$ws.Cells.Item(24,19).Value = "22/08/2023 14:46"
$ws.Cells.Item(24,20).Value = 4.63
$ws.Cells.Item(24,21).Value = "02/09/2023 18:24"
$ws.Cells.Item(24,22).Value = "https://www.betexplorer.com/football/italy/serie-a/bologna-cagliari/r1iTvehJ/"
$ws.Cells.Item(25,6).Value = "Udinese"
$ws.Cells.Item(25,7).Value = 0
$ws.Cells.Item(25,8).Value = "Frosinone"
$ws.Cells.Item(25,9).Value = 0
$ws.Cells.Item(25,10).Value = 1.74
$ws.Cells.Item(25,11).Value = "22/08/2023 14:47"
$ws.Cells.Item(25,12).Value = 1.93
$ws.Cells.Item(25,13).Value = "02/09/2023 18:23"
$ws.Cells.Item(25,14).Value = 3.72
$ws.Cells.Item(25,15).Value = "22/08/2023 14:47"
$ws.Cells.Item(25,16).Value = 3.62
$ws.Cells.Item(25,17).Value = "02/09/2023 18:23"
$ws.Cells.Item(25,18).Value = 4.61
$ws.Cells.Item(25,19).Value = "22/08/2023 14:47"
$ws.Cells.Item(25,20).Value = 4.35
$ws.Cells.Item(25,21).Value = "02/09/2023 18:23"
$ws.Cells.Item(25,22).Value = "https://www.betexplorer.com/football/italy/serie-a/udinese-frosinone/rFW1gf7C/"

# Rows 53, 54, 55
$ws.Cells.Item(53,6).Value = "Empoli"
$ws.Cells.Item(53,7).Value = 1
$ws.Cells.Item(53,8).Value = "Salernitana"
$ws.Cells.Item(53,9).Value = 0
$ws.Cells.Item(53,10).Value = 2.25
$ws.Cells.Item(53,11).Value = "17/09/2023 09:02"
$ws.Cells.Item(53,12).Value = 2.8
$ws.Cells.Item(53,13).Value = "27/09/2023 18:28"
$ws.Cells.Item(53,14).Value = 3.02
$ws.Cells.Item(53,15).Value = "17/09/2023 09:02"
$ws.Cells.Item(53,16).Value = 3.19
$ws.Cells.Item(53,17).Value = "27/09/2023 18:26"
$ws.Cells.Item(53,18).Value = 3.58
$ws.Cells.Item(53,19).Value = "17/09/2023 09:02"
$ws.Cells.Item(53,20).Value = 2.85
$ws.Cells.Item(53,21).Value = "27/09/2023 18:28"
$ws.Cells.Item(53,22).Value = "https://www.betexplorer.com/football/italy/serie-a/empoli-salernitana/jNpkpBod/"
$ws.Cells.Item(54,6).Value = "Verona"
$ws.Cells.Item(54,7).Value = 0
$ws.Cells.Item(54,8).Value = "Atalanta"
$ws.Cells.Item(54,9).Value = 1
$ws.Cells.Item(54,10).Value = 3.96
$ws.Cells.Item(54,11).Value = "17/09/2023 09:02"
$ws.Cells.Item(54,12).Value = 4.48
$ws.Cells.Item(54,13).Value = "27/09/2023 18:29"
$ws.Cells.Item(54,14).Value = 3.76
$ws.Cells.Item(54,15).Value = "17/09/2023 09:02"
$ws.Cells.Item(54,16).Value = 3.66
$ws.Cells.Item(54,17).Value = "27/09/2023 18:29"
$ws.Cells.Item(54,18).Value = 1.85
$ws.Cells.Item(54,19).Value = "17/09/2023 09:02"
$ws.Cells.Item(54,20).Value = 1.9
$ws.Cells.Item(54,21).Value = "27/09/2023 18:29"
$ws.Cells.Item(54,22).Value = "https://www.betexplorer.com/football/italy/serie-a/verona-atalanta/4QAObA0k/"
$ws.Cells.Item(55,6).Value = "Cagliari"
$ws.Cells.Item(55,7).Value = 1
$ws.Cells.Item(55,8).Value = "AC Milan"
$ws.Cells.Item(55,9).Value = 3
$ws.Cells.Item(55,10).Value = 5.55
$ws.Cells.Item(55,11).Value = "23/09/2023 09:29"
$ws.Cells.Item(55,12).Value = 4.73
$ws.Cells.Item(55,13).Value = "27/09/2023 18:10"
$ws.Cells.Item(55,14).Value = 4.14
$ws.Cells.Item(55,15).Value = "23/09/2023 09:29"
$ws.Cells.Item(55,16).Value = 3.43
$ws.Cells.Item(55,17).Value = "27/09/2023 18:10"
$ws.Cells.Item(55,18).Value = 1.56
$ws.Cells.Item(55,19).Value = "23/09/2023 09:29"
$ws.Cells.Item(55,20).Value = 1.92
$ws.Cells.Item(55,21).Value = "27/09/2023 18:10"
$ws.Cells.Item(55,22).Value = "https://www.betexplorer.com/football/italy/serie-a/cagliari-ac-milan/CWoooiWk/"

# Rows 59, 60
$ws.Cells.Item(59,6).Value = "Frosinone"
$ws.Cells.Item(59,7).Value = 1
$ws.Cells.Item(59,8).Value = "Fiorentina"
$ws.Cells.Item(59,9).Value = 1
$ws.Cells.Item(59,10).Value = 4.05
$ws.Cells.Item(59,11).Value = "23/09/2023 09:30"
$ws.Cells.Item(59,12).Value = 3.97
$ws.Cells.Item(59,13).Value = "28/09/2023 18:22"
$ws.Cells.Item(59,14).Value = 3.7
$ws.Cells.Item(59,15).Value = "23/09/2023 09:30"
$ws.Cells.Item(59,16).Value = 3.65
$ws.Cells.Item(59,17).Value = "28/09/2023 18:09"
$ws.Cells.Item(59,18).Value = 1.85
$ws.Cells.Item(59,19).Value = "23/09/2023 09:30"
$ws.Cells.Item(59,20).Value = 2.01
$ws.Cells.Item(59,21).Value = "28/09/2023 18:09"
$ws.Cells.Item(59,22).Value = "https://www.betexplorer.com/football/italy/serie-a/frosinone-fiorentina/4QYW8fOd/"
$ws.Cells.Item(60,6).Value = "Monza"
$ws.Cells.Item(60,7).Value = 0
$ws.Cells.Item(60,8).Value = "Bologna"
$ws.Cells.Item(60,9).Value = 0
$ws.Cells.Item(60,10).Value = 2.43
$ws.Cells.Item(60,11).Value = "17/09/2023 09:02"
$ws.Cells.Item(60,12).Value = 2.31
$ws.Cells.Item(60,13).Value = "28/09/2023 18:29"
$ws.Cells.Item(60,14).Value = 3.36
$ws.Cells.Item(60,15).Value = "17/09/2023 09:02"
$ws.Cells.Item(60,16).Value = 3.45
$ws.Cells.Item(60,17).Value = "28/09/2023 18:25"
$ws.Cells.Item(60,18).Value = 2.89
$ws.Cells.Item(60,19).Value = "17/09/2023 09:02"
$ws.Cells.Item(60,20).Value = 3.35
$ws.Cells.Item(60,21).Value = "28/09/2023 18:29"
$ws.Cells.Item(60,22).Value = "https://www.betexplorer.com/football/italy/serie-a/monza-bologna/n9Ui4WxS/"

# Rows 69, 70
$ws.Cells.Item(69,6).Value = "Sassuolo"
$ws.Cells.Item(69,7).Value = 0
$ws.Cells.Item(69,8).Value = "Monza"
$ws.Cells.Item(69,9).Value = 1
$ws.Cells.Item(69,10).Value = 2.16
$ws.Cells.Item(69,11).Value = "22/09/2023 19:02"
$ws.Cells.Item(69,12).Value = 2.16
$ws.Cells.Item(69,13).Value = "02/10/2023 18:22"
$ws.Cells.Item(69,14).Value = 3.72
$ws.Cells.Item(69,15).Value = "22/09/2023 19:02"
$ws.Cells.Item(69,16).Value = 3.87
$ws.Cells.Item(69,17).Value = "02/10/2023 18:21"
$ws.Cells.Item(69,18).Value = 3.08
$ws.Cells.Item(69,19).Value = "22/09/2023 19:02"
$ws.Cells.Item(69,20).Value = 3.31
$ws.Cells.Item(69,21).Value = "02/10/2023 18:22"
$ws.Cells.Item(69,22).Value = "https://www.betexplorer.com/football/italy/serie-a/sassuolo-monza/ppWcD9Fr/"
$ws.Cells.Item(70,6).Value = "Torino"
$ws.Cells.Item(70,7).Value = 0
$ws.Cells.Item(70,8).Value = "Verona"
$ws.Cells.Item(70,9).Value = 0
$ws.Cells.Item(70,10).Value = 1.67
$ws.Cells.Item(70,11).Value = "22/09/2023 19:02"
$ws.Cells.Item(70,12).Value = 1.76
$ws.Cells.Item(70,13).Value = "02/10/2023 18:23"
$ws.Cells.Item(70,14).Value = 3.69
$ws.Cells.Item(70,15).Value = "22/09/2023 19:02"
$ws.Cells.Item(70,16).Value = 3.47
$ws.Cells.Item(70,17).Value = "02/10/2023 18:23"
$ws.Cells.Item(70,18).Value = 5.24
$ws.Cells.Item(70,19).Value = "22/09/2023 19:02"
$ws.Cells.Item(70,20).Value = 5.81
$ws.Cells.Item(70,21).Value = "02/10/2023 18:29"
$ws.Cells.Item(70,22).Value = "https://www.betexplorer.com/football/italy/serie-a/torino-verona/hAV1CTUl/"

# Rows 86, 87
$ws.Cells.Item(86,6).Value = "Salernitana"
$ws.Cells.Item(86,7).Value = 2
$ws.Cells.Item(86,8).Value = "Cagliari"
$ws.Cells.Item(86,9).Value = 2
$ws.Cells.Item(86,10).Value = 2.47
$ws.Cells.Item(86,11).Value = "06/10/2023 09:04"
$ws.Cells.Item(86,12).Value = 2.34
$ws.Cells.Item(86,13).Value = "22/10/2023 14:59"
$ws.Cells.Item(86,14).Value = 3.27
$ws.Cells.Item(86,15).Value = "06/10/2023 09:04"
$ws.Cells.Item(86,16).Value = 3.38
$ws.Cells.Item(86,17).Value = "22/10/2023 14:59"
$ws.Cells.Item(86,18).Value = 2.91
$ws.Cells.Item(86,19).Value = "06/10/2023 09:04"
$ws.Cells.Item(86,20).Value = 3.35
$ws.Cells.Item(86,21).Value = "22/10/2023 14:59"
$ws.Cells.Item(86,22).Value = "https://www.betexplorer.com/football/italy/serie-a/salernitana-cagliari/CYUDvnEs/"
$ws.Cells.Item(87,6).Value = "Bologna"
$ws.Cells.Item(87,7).Value = 2
$ws.Cells.Item(87,8).Value = "Frosinone"
$ws.Cells.Item(87,9).Value = 1
$ws.Cells.Item(87,10).Value = 1.67
$ws.Cells.Item(87,11).Value = "06/10/2023 09:04"
$ws.Cells.Item(87,12).Value = 1.74
$ws.Cells.Item(87,13).Value = "22/10/2023 14:59"
$ws.Cells.Item(87,14).Value = 3.9
$ws.Cells.Item(87,15).Value = "06/10/2023 09:04"
$ws.Cells.Item(87,16).Value = 3.92
$ws.Cells.Item(87,17).Value = "22/10/2023 14:59"
$ws.Cells.Item(87,18).Value = 5.49
$ws.Cells.Item(87,19).Value = "06/10/2023 09:04"
$ws.Cells.Item(87,20).Value = 5.1
$ws.Cells.Item(87,21).Value = "22/10/2023 14:59"
$ws.Cells.Item(87,22).Value = "https://www.betexplorer.com/football/italy/serie-a/bologna-frosinone/0f1eM4bD/"

# ---- Part 2: append the 10 new matches scraped on 31-10-2023 (rows 92-101) ----

# Copy row 91's formatting (bold/bordered index column, datetime style on data_partida)
# down across the new rows before writing their values.
$ws.Range("A91:V91").Copy()
$ws.Range("A92:V101").PasteSpecial(-4122)

# Row 92 (Indice 91)
$ws.Cells.Item(92,1).Value = 91
$ws.Cells.Item(92,2).Value = "italy"
$ws.Cells.Item(92,3).Value = "serie-a"
$ws.Cells.Item(92,4).Value = "2023-2024"
$ws.Cells.Item(92,5).Value = 45226.86458333334
$ws.Cells.Item(92,6).Value = "Genoa"
$ws.Cells.Item(92,7).Value = 1
$ws.Cells.Item(92,8).Value = "Salernitana"
$ws.Cells.Item(92,9).Value = 0
$ws.Cells.Item(92,10).Value = 1.61
$ws.Cells.Item(92,11).Value = "11/10/2023 14:09"
$ws.Cells.Item(92,12).Value = 1.75
$ws.Cells.Item(92,13).Value = "27/10/2023 20:25"
$ws.Cells.Item(92,14).Value = 3.96
$ws.Cells.Item(92,15).Value = "11/10/2023 14:09"
$ws.Cells.Item(92,16).Value = 3.69
$ws.Cells.Item(92,17).Value = "27/10/2023 20:41"
$ws.Cells.Item(92,18).Value = 5.28
$ws.Cells.Item(92,19).Value = "11/10/2023 14:09"
$ws.Cells.Item(92,20).Value = 5.29
$ws.Cells.Item(92,21).Value = "27/10/2023 20:41"
$ws.Cells.Item(92,22).Value = "https://www.betexplorer.com/football/italy/serie-a/genoa-salernitana/n5fwsuzt/"

# Row 93 (Indice 92)
$ws.Cells.Item(93,1).Value = 92
$ws.Cells.Item(93,2).Value = "italy"
$ws.Cells.Item(93,3).Value = "serie-a"
$ws.Cells.Item(93,4).Value = "2023-2024"
$ws.Cells.Item(93,5).Value = 45227.625
$ws.Cells.Item(93,6).Value = "Sassuolo"
$ws.Cells.Item(93,7).Value = 1
$ws.Cells.Item(93,8).Value = "Bologna"
$ws.Cells.Item(93,9).Value = 1
$ws.Cells.Item(93,10).Value = 2.4
$ws.Cells.Item(93,11).Value = "10/10/2023 14:02"
$ws.Cells.Item(93,12).Value = 2.58
$ws.Cells.Item(93,13).Value = "28/10/2023 14:58"
$ws.Cells.Item(93,14).Value = 3.51
$ws.Cells.Item(93,15).Value = "10/10/2023 14:02"
$ws.Cells.Item(93,16).Value = 3.46
$ws.Cells.Item(93,17).Value = "28/10/2023 14:58"
$ws.Cells.Item(93,18).Value = 3.01
$ws.Cells.Item(93,19).Value = "10/10/2023 14:02"
$ws.Cells.Item(93,20).Value = 2.85
$ws.Cells.Item(93,21).Value = "28/10/2023 14:58"
$ws.Cells.Item(93,22).Value = "https://www.betexplorer.com/football/italy/serie-a/sassuolo-bologna/lzv7zeJN/"

# Row 94 (Indice 93)
$ws.Cells.Item(94,1).Value = 93
$ws.Cells.Item(94,2).Value = "italy"
$ws.Cells.Item(94,3).Value = "serie-a"
$ws.Cells.Item(94,4).Value = "2023-2024"
$ws.Cells.Item(94,5).Value = 45227.75
$ws.Cells.Item(94,6).Value = "Lecce"
$ws.Cells.Item(94,7).Value = 0
$ws.Cells.Item(94,8).Value = "Torino"
$ws.Cells.Item(94,9).Value = 1
$ws.Cells.Item(94,10).Value = 2.99
$ws.Cells.Item(94,11).Value = "10/10/2023 14:02"
$ws.Cells.Item(94,12).Value = 2.79
$ws.Cells.Item(94,13).Value = "28/10/2023 17:55"
$ws.Cells.Item(94,14).Value = 2.95
$ws.Cells.Item(94,15).Value = "10/10/2023 14:02"
$ws.Cells.Item(94,16).Value = 2.93
$ws.Cells.Item(94,17).Value = "28/10/2023 17:55"
$ws.Cells.Item(94,18).Value = 2.77
$ws.Cells.Item(94,19).Value = "10/10/2023 14:02"
$ws.Cells.Item(94,20).Value = 3.05
$ws.Cells.Item(94,21).Value = "28/10/2023 17:58"
$ws.Cells.Item(94,22).Value = "https://www.betexplorer.com/football/italy/serie-a/lecce-torino/EHofwcZ4/"

# Row 95 (Indice 94)
$ws.Cells.Item(95,1).Value = 94
$ws.Cells.Item(95,2).Value = "italy"
$ws.Cells.Item(95,3).Value = "serie-a"
$ws.Cells.Item(95,4).Value = "2023-2024"
$ws.Cells.Item(95,5).Value = 45227.86458333334
$ws.Cells.Item(95,6).Value = "Juventus"
$ws.Cells.Item(95,7).Value = 1
$ws.Cells.Item(95,8).Value = "Verona"
$ws.Cells.Item(95,9).Value = 0
$ws.Cells.Item(95,10).Value = 1.43
$ws.Cells.Item(95,11).Value = "10/10/2023 14:02"
$ws.Cells.Item(95,12).Value = 1.41
$ws.Cells.Item(95,13).Value = "28/10/2023 20:43"
$ws.Cells.Item(95,14).Value = 4.67
$ws.Cells.Item(95,15).Value = "10/10/2023 14:02"
$ws.Cells.Item(95,16).Value = 4.84
$ws.Cells.Item(95,17).Value = "28/10/2023 20:44"
$ws.Cells.Item(95,18).Value = 8.08
$ws.Cells.Item(95,19).Value = "10/10/2023 14:02"
$ws.Cells.Item(95,20).Value = 8.69
$ws.Cells.Item(95,21).Value = "28/10/2023 20:44"
$ws.Cells.Item(95,22).Value = "https://www.betexplorer.com/football/italy/serie-a/juventus-verona/8KrnuJ4h/"

# Row 96 (Indice 95)
$ws.Cells.Item(96,1).Value = 95
$ws.Cells.Item(96,2).Value = "italy"
$ws.Cells.Item(96,3).Value = "serie-a"
$ws.Cells.Item(96,4).Value = "2023-2024"
$ws.Cells.Item(96,5).Value = 45228.52083333334
$ws.Cells.Item(96,6).Value = "Cagliari"
$ws.Cells.Item(96,7).Value = 4
$ws.Cells.Item(96,8).Value = "Frosinone"
$ws.Cells.Item(96,9).Value = 3
$ws.Cells.Item(96,10).Value = 2.4
$ws.Cells.Item(96,11).Value = "11/10/2023 14:09"
$ws.Cells.Item(96,12).Value = 2.19
$ws.Cells.Item(96,13).Value = "29/10/2023 12:26"
$ws.Cells.Item(96,14).Value = 3.36
$ws.Cells.Item(96,15).Value = "11/10/2023 14:09"
$ws.Cells.Item(96,16).Value = 3.46
$ws.Cells.Item(96,17).Value = "29/10/2023 12:24"
$ws.Cells.Item(96,18).Value = 3.14
$ws.Cells.Item(96,19).Value = "11/10/2023 14:09"
$ws.Cells.Item(96,20).Value = 3.54
$ws.Cells.Item(96,21).Value = "29/10/2023 12:28"
$ws.Cells.Item(96,22).Value = "https://www.betexplorer.com/football/italy/serie-a/cagliari-frosinone/MPwSZNSC/"

# Row 97 (Indice 96)
$ws.Cells.Item(97,1).Value = 96
$ws.Cells.Item(97,2).Value = "italy"
$ws.Cells.Item(97,3).Value = "serie-a"
$ws.Cells.Item(97,4).Value = "2023-2024"
$ws.Cells.Item(97,5).Value = 45228.625
$ws.Cells.Item(97,6).Value = "Monza"
$ws.Cells.Item(97,7).Value = 1
$ws.Cells.Item(97,8).Value = "Udinese"
$ws.Cells.Item(97,9).Value = 1
$ws.Cells.Item(97,10).Value = 2.04
$ws.Cells.Item(97,11).Value = "10/10/2023 14:02"
$ws.Cells.Item(97,12).Value = 2.13
$ws.Cells.Item(97,13).Value = "29/10/2023 14:57"
$ws.Cells.Item(97,14).Value = 3.56
$ws.Cells.Item(97,15).Value = "10/10/2023 14:02"
$ws.Cells.Item(97,16).Value = 3.53
$ws.Cells.Item(97,17).Value = "29/10/2023 14:57"
$ws.Cells.Item(97,18).Value = 3.51
$ws.Cells.Item(97,19).Value = "10/10/2023 14:02"
$ws.Cells.Item(97,20).Value = 3.64
$ws.Cells.Item(97,21).Value = "29/10/2023 14:57"
$ws.Cells.Item(97,22).Value = "https://www.betexplorer.com/football/italy/serie-a/monza-udinese/h8pbxHlB/"

# Row 98 (Indice 97)
$ws.Cells.Item(98,1).Value = 97
$ws.Cells.Item(98,2).Value = "italy"
$ws.Cells.Item(98,3).Value = "serie-a"
$ws.Cells.Item(98,4).Value = "2023-2024"
$ws.Cells.Item(98,5).Value = 45228.75
$ws.Cells.Item(98,6).Value = "Inter"
$ws.Cells.Item(98,7).Value = 1
$ws.Cells.Item(98,8).Value = "AS Roma"
$ws.Cells.Item(98,9).Value = 0
$ws.Cells.Item(98,10).Value = 1.74
$ws.Cells.Item(98,11).Value = "10/10/2023 14:02"
$ws.Cells.Item(98,12).Value = 1.54
$ws.Cells.Item(98,13).Value = "29/10/2023 17:52"
$ws.Cells.Item(98,14).Value = 3.82
$ws.Cells.Item(98,15).Value = "10/10/2023 14:02"
$ws.Cells.Item(98,16).Value = 4.21
$ws.Cells.Item(98,17).Value = "29/10/2023 17:58"
$ws.Cells.Item(98,18).Value = 5.02
$ws.Cells.Item(98,19).Value = "10/10/2023 14:02"
$ws.Cells.Item(98,20).Value = 6.83
$ws.Cells.Item(98,21).Value = "29/10/2023 17:58"
$ws.Cells.Item(98,22).Value = "https://www.betexplorer.com/football/italy/serie-a/inter-as-roma/Qegstakn/"

# Row 99 (Indice 98)
$ws.Cells.Item(99,1).Value = 98
$ws.Cells.Item(99,2).Value = "italy"
$ws.Cells.Item(99,3).Value = "serie-a"
$ws.Cells.Item(99,4).Value = "2023-2024"
$ws.Cells.Item(99,5).Value = 45228.86458333334
$ws.Cells.Item(99,6).Value = "Napoli"
$ws.Cells.Item(99,7).Value = 2
$ws.Cells.Item(99,8).Value = "AC Milan"
$ws.Cells.Item(99,9).Value = 2
$ws.Cells.Item(99,10).Value = 1.97
$ws.Cells.Item(99,11).Value = "10/10/2023 14:02"
$ws.Cells.Item(99,12).Value = 2.27
$ws.Cells.Item(99,13).Value = "29/10/2023 20:44"
$ws.Cells.Item(99,14).Value = 3.59
$ws.Cells.Item(99,15).Value = "10/10/2023 14:02"
$ws.Cells.Item(99,16).Value = 3.43
$ws.Cells.Item(99,17).Value = "29/10/2023 20:40"
$ws.Cells.Item(99,18).Value = 4.01
$ws.Cells.Item(99,19).Value = "10/10/2023 14:02"
$ws.Cells.Item(99,20).Value = 3.39
$ws.Cells.Item(99,21).Value = "29/10/2023 20:44"
$ws.Cells.Item(99,22).Value = "https://www.betexplorer.com/football/italy/serie-a/napoli-ac-milan/Aqu3yy4H/"

# Row 100 (Indice 99)
$ws.Cells.Item(100,1).Value = 99
$ws.Cells.Item(100,2).Value = "italy"
$ws.Cells.Item(100,3).Value = "serie-a"
$ws.Cells.Item(100,4).Value = "2023-2024"
$ws.Cells.Item(100,5).Value = 45229.77083333334
$ws.Cells.Item(100,6).Value = "Empoli"
$ws.Cells.Item(100,7).Value = 0
$ws.Cells.Item(100,8).Value = "Atalanta"
$ws.Cells.Item(100,9).Value = 3
$ws.Cells.Item(100,10).Value = 5
$ws.Cells.Item(100,11).Value = "10/10/2023 14:02"
$ws.Cells.Item(100,12).Value = 5.06
$ws.Cells.Item(100,13).Value = "30/10/2023 18:29"
$ws.Cells.Item(100,14).Value = 4.17
$ws.Cells.Item(100,15).Value = "10/10/2023 14:02"
$ws.Cells.Item(100,16).Value = 4.15
$ws.Cells.Item(100,17).Value = "30/10/2023 18:28"
$ws.Cells.Item(100,18).Value = 1.61
$ws.Cells.Item(100,19).Value = "10/10/2023 14:02"
$ws.Cells.Item(100,20).Value = 1.67
$ws.Cells.Item(100,21).Value = "30/10/2023 18:11"
$ws.Cells.Item(100,22).Value = "https://www.betexplorer.com/football/italy/serie-a/empoli-atalanta/0GxWYsrJ/"

# Row 101 (Indice 100)
$ws.Cells.Item(101,1).Value = 100
$ws.Cells.Item(101,2).Value = "italy"
$ws.Cells.Item(101,3).Value = "serie-a"
$ws.Cells.Item(101,4).Value = "2023-2024"
$ws.Cells.Item(101,5).Value = 45229.86458333334
$ws.Cells.Item(101,6).Value = "Lazio"
$ws.Cells.Item(101,7).Value = 1
$ws.Cells.Item(101,8).Value = "Fiorentina"
$ws.Cells.Item(101,9).Value = 0
$ws.Cells.Item(101,10).Value = 2.17
$ws.Cells.Item(101,11).Value = "10/10/2023 14:02"
$ws.Cells.Item(101,12).Value = 2.33
$ws.Cells.Item(101,13).Value = "30/10/2023 20:42"
$ws.Cells.Item(101,14).Value = 3.51
$ws.Cells.Item(101,15).Value = "10/10/2023 14:02"
$ws.Cells.Item(101,16).Value = 3.38
$ws.Cells.Item(101,17).Value = "30/10/2023 20:43"
$ws.Cells.Item(101,18).Value = 3.24
$ws.Cells.Item(101,19).Value = "10/10/2023 14:02"
$ws.Cells.Item(101,20).Value = 3.3
$ws.Cells.Item(101,21).Value = "30/10/2023 20:44"
$ws.Cells.Item(101,22).Value = "https://www.betexplorer.com/football/italy/serie-a/lazio-fiorentina/KfnjvwKb/"

